$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1184.0625
$ws.Range("I19").Value = 823.2857
$ws.Range("J19").Value = 1464.6666
$ws.Range("K19").Value = 823.2857
$ws.Range("L19").Value = 1464.6666
$ws.Range("M19").Value = -648.2857
$ws.Range("N19").Value = -1814.6666

$ws.Range("H33").Value = 354.17648
$ws.Range("I33").Value = 366.9375
$ws.Range("K33").Value = 366.9375
$ws.Range("M33").Value = -137.9375

$ws.Range("H53").Value = 1370.5834
$ws.Range("I53").Value = 1030.6
$ws.Range("J53").Value = 1613.4286
$ws.Range("K53").Value = 1030.6
$ws.Range("L53").Value = 1613.4286
$ws.Range("M53").Value = -393.5999999999999
$ws.Range("N53").Value = -2887.4286

$ws.Range("H62").Value = 5099.4546
$ws.Range("I62").Value = 4772.2354
$ws.Range("K62").Value = 4772.2354
$ws.Range("M62").Value = -4148.2354

$ws.Range("H65").Value = 5099.4546
$ws.Range("I65").Value = 4772.2354
$ws.Range("K65").Value = 23861.177
$ws.Range("M65").Value = -20741.177

$ws.Range("H98").Value = 567.2727
$ws.Range("I98").Value = 630.125
$ws.Range("K98").Value = 630.125
$ws.Range("M98").Value = 867.875

$ws.Range("H122").Value = 567.2727
$ws.Range("I122").Value = 630.125
$ws.Range("K122").Value = 1890.375
$ws.Range("M122").Value = 559.625

$ws.Range("H138").Value = 3132.2683
$ws.Range("I138").Value = 1969.2727
$ws.Range("J138").Value = 3558.7
$ws.Range("K138").Value = 5907.8181
$ws.Range("L138").Value = 10676.1
$ws.Range("M138").Value = -767.8181000000004
$ws.Range("N138").Value = -20956.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8895.093000000001
$ws.Range("I32").Value = 8219.75
$ws.Range("K32").Value = 8219.75
$ws.Range("M32").Value = -7932.75

$ws.Range("H45").Value = 3569.25
$ws.Range("J45").Value = 4259
$ws.Range("L45").Value = 4259
$ws.Range("N45").Value = -5013

$ws.Range("H60").Value = 89522.5
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()

$ws.Range("H61").Value = 7504874
$ws.Range("I61").Value = 10005285
$ws.Range("J61").Value = 1253847.5
$ws.Range("K61").Value = 10005285
$ws.Range("L61").Value = 1253847.5
$ws.Range("M61").Value = -10005073
$ws.Range("N61").Value = -1254271.5

$ws.Range("H102").Value = 8991.111000000001
$ws.Range("I102").Value = 8124.75
$ws.Range("J102").Value = 9684.200000000001
$ws.Range("K102").Value = 8124.75
$ws.Range("L102").Value = 9684.200000000001
$ws.Range("M102").Value = -6502.75
$ws.Range("N102").Value = -12928.2

$ws.Range("H122").Value = 4688.385
$ws.Range("I122").Value = 4787.5
$ws.Range("J122").Value = 3499
$ws.Range("K122").Value = 14362.5
$ws.Range("L122").Value = 10497
$ws.Range("M122").Value = -11912.5
$ws.Range("N122").Value = -15397

$ws.Range("H136").Value = 7504874
$ws.Range("I136").Value = 10005285
$ws.Range("J136").Value = 1253847.5
$ws.Range("K136").Value = 30015855
$ws.Range("L136").Value = 3761542.5
$ws.Range("M136").Value = -30013305
$ws.Range("N136").Value = -3766642.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2054.2273
$ws.Range("I94").Value = 2368.875
$ws.Range("K94").Value = 2368.875
$ws.Range("M94").Value = -1917.875

$ws.Range("H107").Value = 2449.0435
$ws.Range("I107").Value = 2767.6667
$ws.Range("K107").Value = 2767.6667
$ws.Range("M107").Value = -847.6667000000002

$ws.Range("H134").Value = 2565919.2
$ws.Range("I134").Value = 1760.1666
$ws.Range("K134").Value = 5280.4998
$ws.Range("M134").Value = -2745.4998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 40002948
$ws.Range("I31").Value = 52634028
$ws.Range("K31").Value = 52634028
$ws.Range("M31").Value = -52633733

$ws.Range("H34").Value = 40002948
$ws.Range("I34").Value = 52634028
$ws.Range("K34").Value = 52634028
$ws.Range("M34").Value = -52633826

$ws.Range("H58").Value = 3840.318
$ws.Range("I58").Value = 3467.375
$ws.Range("K58").Value = 3467.375
$ws.Range("M58").Value = -3264.375

$ws.Range("H117").Value = 99999
$ws.Range("J117").Value = 99999
$ws.Range("L117").Value = 99999
$ws.Range("N117").Value = -109177

$ws.Range("H122").Value = 2063.3076
$ws.Range("I122").Value = 2080.9
$ws.Range("K122").Value = 6242.700000000001
$ws.Range("M122").Value = -3792.700000000001

$ws.Range("H136").Value = 3840.318
$ws.Range("I136").Value = 3467.375
$ws.Range("K136").Value = 10402.125
$ws.Range("M136").Value = -7852.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 10326.667
$ws.Range("I3").Value = 5725.4
$ws.Range("J3").Value = 33333
$ws.Range("K3").Value = 17176.2
$ws.Range("L3").Value = 99999
$ws.Range("M3").Value = -17064.2
$ws.Range("N3").Value = -100223

$ws.Range("H5").Value = 505.6
$ws.Range("I5").Value = 346.1875
$ws.Range("J5").Value = 1143.25
$ws.Range("K5").Value = 1038.5625
$ws.Range("L5").Value = 3429.75
$ws.Range("M5").Value = -926.5625
$ws.Range("N5").Value = -3653.75

$ws.Range("H49").Value = 12111
$ws.Range("I49").Value = 1500
$ws.Range("K49").Value = 4500
$ws.Range("M49").Value = -4344

$ws.Range("H112").Value = 21111
$ws.Range("I112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("M112").ClearContents()

$ws.Range("H135").Value = 505.6
$ws.Range("I135").Value = 346.1875
$ws.Range("J135").Value = 1143.25
$ws.Range("K135").Value = 3115.6875
$ws.Range("L135").Value = 10289.25
$ws.Range("M135").Value = -580.6875
$ws.Range("N135").Value = -15359.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11177.38
$ws.Range("J70").Value = 11228.167
$ws.Range("L70").Value = 11228.167
$ws.Range("N70").Value = -11768.167

$ws.Range("H73").Value = 11177.38
$ws.Range("J73").Value = 11228.167
$ws.Range("L73").Value = 11228.167
$ws.Range("N73").Value = -13100.167

$ws.Range("H102").Value = 1499.5312
$ws.Range("I102").Value = 1519.5667
$ws.Range("K102").Value = 1519.5667
$ws.Range("M102").Value = 102.4332999999999

$ws.Range("H122").Value = 8745.235000000001
$ws.Range("I122").Value = 7900
$ws.Range("J122").Value = 9336.9
$ws.Range("K122").Value = 23700
$ws.Range("L122").Value = 28010.7
$ws.Range("M122").Value = -21250
$ws.Range("N122").Value = -32910.7

$ws.Range("H126").Value = 5208.385
$ws.Range("I126").Value = 4620.909
$ws.Range("K126").Value = 13862.727
$ws.Range("M126").Value = -11392.727

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7896.316
$ws.Range("I7").Value = 7648.9414
$ws.Range("K7").Value = 7648.9414
$ws.Range("M7").Value = -7536.9414

$ws.Range("H16").Value = 1906.2693
$ws.Range("I16").Value = 1888.3158
$ws.Range("J16").Value = 1955
$ws.Range("K16").Value = 1888.3158
$ws.Range("L16").Value = 1955
$ws.Range("M16").Value = -1718.3158
$ws.Range("N16").Value = -2295

$ws.Range("H40").Value = 7200.9287
$ws.Range("I40").Value = 7234.4165
$ws.Range("J40").Value = 7000
$ws.Range("K40").Value = 7234.4165
$ws.Range("L40").Value = 7000
$ws.Range("M40").Value = -7098.4165
$ws.Range("N40").Value = -7272

$ws.Range("H57").Value = 27114.111
$ws.Range("I57").Value = 27114.111
$ws.Range("K57").Value = 27114.111
$ws.Range("M57").Value = -26548.111

$ws.Range("H100").Value = 22754264
$ws.Range("I100").Value = 2908.5
$ws.Range("J100").Value = 50055890
$ws.Range("K100").Value = 2908.5
$ws.Range("L100").Value = 50055890
$ws.Range("M100").Value = -2367.5
$ws.Range("N100").Value = -50056972

$ws.Range("H122").Value = 4664.4893
$ws.Range("I122").Value = 3483.5134
$ws.Range("K122").Value = 10450.5402
$ws.Range("M122").Value = -8000.540199999999

$ws.Range("H126").Value = 7896.316
$ws.Range("I126").Value = 7648.9414
$ws.Range("K126").Value = 22946.8242
$ws.Range("M126").Value = -20476.8242

$ws.Range("H136").Value = 5472.1177
$ws.Range("I136").Value = 2103.2
$ws.Range("J136").Value = 10284.857
$ws.Range("K136").Value = 6309.599999999999
$ws.Range("L136").Value = 30854.571
$ws.Range("M136").Value = -3759.599999999999
$ws.Range("N136").Value = -35954.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 12266
$ws.Range("J62").Value = 17150
$ws.Range("L62").Value = 17150
$ws.Range("N62").Value = -18398

$ws.Range("H65").Value = 12266
$ws.Range("J65").Value = 17150
$ws.Range("L65").Value = 85750
$ws.Range("N65").Value = -91990

$ws.Range("H81").Value = 2150.2727
$ws.Range("I81").Value = 2217.4
$ws.Range("K81").Value = 4434.8
$ws.Range("M81").Value = -3373.8

$ws.Range("H84").Value = 2150.2727
$ws.Range("I84").Value = 2217.4
$ws.Range("K84").Value = 22174
$ws.Range("M84").Value = -16870

$ws.Range("H126").Value = 3831.0435
$ws.Range("I126").Value = 4579.357
$ws.Range("J126").Value = 2667
$ws.Range("K126").Value = 13738.071
$ws.Range("L126").Value = 8001
$ws.Range("M126").Value = -11268.071
$ws.Range("N126").Value = -12941

$ws.Range("H132").Value = 296940.3
$ws.Range("I132").Value = 2997.3928
$ws.Range("J132").Value = 1668674
$ws.Range("K132").Value = 8992.178400000001
$ws.Range("L132").Value = 5006022
$ws.Range("M132").Value = -6462.178400000001
$ws.Range("N132").Value = -5011082

$ws.Range("H136").Value = 419012.22
$ws.Range("I136").Value = 2164.9
$ws.Range("K136").Value = 6494.700000000001
$ws.Range("M136").Value = -3944.700000000001
